$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name stays the internal sheet name "RGossF-HW50.xpc" -> "RGossF")
$ws.Name = "RGossF"

# Add a new data row (row 16), mirroring the structure/content of row 15
# Copy formatting (font/border/alignment) from A15 so A16 matches the style
# used for the other "index" column cells.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1
